# Update NATMI LR-pair sheet with newly recomputed TPM-based expression values.
#
# Only the "raw" per-cluster ligand/receptor average & total expression values
# (columns G, H for the sending/ligand cluster and M, N for the target/receptor
# cluster) actually come from new data. All of the other numeric columns on the
# sheet (I, J, O, P, Q, R, S, T) are specificity / edge-weight figures that are
# purely derived from those raw values, so they are recomputed here rather than
# hard-coded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand (sending cluster) average/total expression values, keyed by the
# "Sending cluster" name (column A).
$newLigand = @{
    "ECs"               = @{ Avg = 3.062550666666667;  Total = 9.187652 }
    "FAPs"              = @{ Avg = 17.74214666666667;  Total = 53.22644 }
    "Inflammatory-Mac"  = @{ Avg = 13.27534766666667;  Total = 39.826043 }
    "MuSCs"             = @{ Avg = 3.455866;            Total = 10.367598 }
    "Resolving-Mac"     = @{ Avg = 7.175465666666668;  Total = 21.526397 }
}

# New receptor (target cluster) average/total expression values, keyed by the
# "Target cluster" name (column D).
$newReceptor = @{
    "ECs"               = @{ Avg = 60.538204;           Total = 181.614612 }
    "FAPs"              = @{ Avg = 10.84181733333333;  Total = 32.525452 }
    "Inflammatory-Mac"  = @{ Avg = 8.850437666666666;  Total = 26.551313 }
    "MuSCs"             = @{ Avg = 1.757142;             Total = 5.271426 }
    "Resolving-Mac"     = @{ Avg = 16.87263033333333;  Total = 50.617891 }
}

$firstRow = 2
$lastRow = 26

# --- Step 1: write the new raw ligand/receptor expression values -----------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sender = $ws.Range("A$r").Value2
    $target = $ws.Range("D$r").Value2

    $ws.Range("G$r").Value = $newLigand[$sender].Avg
    $ws.Range("H$r").Value = $newLigand[$sender].Total

    $ws.Range("M$r").Value = $newReceptor[$target].Avg
    $ws.Range("N$r").Value = $newReceptor[$target].Total
}

# --- Step 2: recompute ligand derived specificity (I, J) -------------------
# specificity = value for this sending cluster / sum of that value over all
# (distinct) sending clusters.
$sumLigandAvg = 0.0
$sumLigandTotal = 0.0
foreach ($key in $newLigand.Keys) {
    $sumLigandAvg += $newLigand[$key].Avg
    $sumLigandTotal += $newLigand[$key].Total
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sender = $ws.Range("A$r").Value2
    $ws.Range("I$r").Value = $newLigand[$sender].Avg / $sumLigandAvg
    $ws.Range("J$r").Value = $newLigand[$sender].Total / $sumLigandTotal
}

# --- Step 3: recompute receptor derived specificity (O, P) -----------------
$sumReceptorAvg = 0.0
$sumReceptorTotal = 0.0
foreach ($key in $newReceptor.Keys) {
    $sumReceptorAvg += $newReceptor[$key].Avg
    $sumReceptorTotal += $newReceptor[$key].Total
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $target = $ws.Range("D$r").Value2
    $ws.Range("O$r").Value = $newReceptor[$target].Avg / $sumReceptorAvg
    $ws.Range("P$r").Value = $newReceptor[$target].Total / $sumReceptorTotal
}

# --- Step 4: recompute edge expression weights (Q, R) -----------------------
# edge average weight = ligand average * receptor average
# edge total weight   = ligand total   * receptor total
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $g = $ws.Range("G$r").Value2
    $h = $ws.Range("H$r").Value2
    $m = $ws.Range("M$r").Value2
    $n = $ws.Range("N$r").Value2

    $ws.Range("Q$r").Value = $g * $m
    $ws.Range("R$r").Value = $h * $n
}

# --- Step 5: recompute edge derived specificity (S, T) ---------------------
# specificity = this edge's weight / sum of that weight across every edge
# (every row) in the sheet.
$sumQ = 0.0
$sumR = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sumQ += $ws.Range("Q$r").Value2
    $sumR += $ws.Range("R$r").Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("S$r").Value = $ws.Range("Q$r").Value2 / $sumQ
    $ws.Range("T$r").Value = $ws.Range("R$r").Value2 / $sumR
}

Write-Output "Updated TPM-derived values for rows $firstRow..$lastRow"
